$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 184 (shifts existing rows 184-238 down to 185-239)
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new data record
$ws.Cells.Item(184, 1).Value = 6
$ws.Cells.Item(184, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(184, 3).Value = "Metropolitana"
$ws.Cells.Item(184, 4).Value = 44782
$ws.Cells.Item(184, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(184, 5).Value = 13
$ws.Cells.Item(184, 6).Value = 100112022
$ws.Cells.Item(184, 7).Value = "Arveja Verde"
$ws.Cells.Item(184, 8).Value = "Perfection"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 220
$ws.Cells.Item(184, 11).Value = 37000
$ws.Cells.Item(184, 12).Value = 38000
$ws.Cells.Item(184, 13).Value = 37455
$ws.Cells.Item(184, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(184, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(184, 16).Value = 1498
$ws.Cells.Item(184, 17).Value = 25
$ws.Cells.Item(184, 18).Value = "Hortaliza"
